# Insert a new data row at row 711 (2026/01/25, 日, 16, 201) into Sheet1.
# All rows from the old 711 through 752 shift down by one (to 712..753),
# which Rows().Insert() handles automatically since it shifts existing
# cells down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 711 (and everything below it) down by one row.
$ws.Rows(711).Insert()

# Fill in the newly inserted row. Column A holds a date-shaped string
# ("2026/01/25") that must stay literal text (matching every other date
# cell in the sheet), so force Text format before assigning it, then
# clear the format back off so no stray style index is left behind.
$ws.Range("A711").NumberFormat = "@"
$ws.Range("A711").Value = "2026/01/25"
$ws.Range("A711").ClearFormats()

$ws.Range("B711").Value = "日"
$ws.Range("C711").Value = 16
$ws.Range("D711").Value = 201
